# This script records a new sale in the cafe point-of-sale workbook:
#   - One "Cafe Helado" is sold, its price increased 6000 -> 8000 and
#     its stock decremented 197 -> 196 on the Productos sheet.
#   - A new row is appended to the Ventas (sales) sheet for the sale.
#   - A new row is appended to the DetalleVentas (sale line items) sheet
#     describing the single line item of the sale.

$wb = $excel.ActiveWorkbook

$saleId    = "c82d54ad-d298-4186-b1ab-20aa106d1339"
$productId = "2e23a440-e5e1-4ce2-a407-c15bad1fd402"
$timestamp = "2025-09-04T21:18:57.812Z"
$saleTime  = "2025-09-04T21:18:57.790Z"
$productName = "Café Helado"
$newPrice  = 8000
$newStock  = 196

# --- Productos: update price, stock and updatedAt for "Cafe Helado" (row 2) ---
$productos = $wb.Worksheets.Item("Productos")
$productos.Range("E2").Value = $newPrice
$productos.Range("F2").Value = $newStock
$productos.Range("I2").Value = $timestamp

# --- Ventas: append the new sale as row 6 ---
$ventas = $wb.Worksheets.Item("Ventas")
$ventas.Range("A6").Value = $saleId
$ventas.Range("C6").Value = "Cliente General"
$ventas.Range("D6").Value = $newPrice
$ventas.Range("E6").Value = $saleTime
$ventas.Range("F6").Value = $saleTime

# --- DetalleVentas: append the sale's single line item as row 6 ---
$detalle = $wb.Worksheets.Item("DetalleVentas")
$detalle.Range("A6").Value = $saleId + "_" + $productId
$detalle.Range("B6").Value = $saleId
$detalle.Range("C6").Value = $productId
$detalle.Range("D6").Value = $productName
$detalle.Range("E6").Value = 1
$detalle.Range("F6").Value = $newPrice
$detalle.Range("G6").Value = $newPrice
